# Auto-generated edit script applying the cryptos.xlsx diff
# (GitHub Actions crypto price refresh, Tue Nov 28 05:55:41 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Preserve the new value as literal text (matches the original
    # inline-string cell), even when it looks like a plain number
    # (e.g. "224.91", "1.00", "2.76"), and leave the cell style index
    # exactly as it was (no s="..." attribute) once done.
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '36.944.56'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '2.001.19'
$ws.Range('E3').Value = '  -2.34%  '
$ws.Range('E4').Value = '  +0.30%  '
Set-TextValue 'D5' '224.91'
$ws.Range('E5').Value = '  -1.94%  '
$ws.Range('E6').Value = '  -1.55%  '
$ws.Range('E7').Value = '  +0.07%  '
Set-TextValue 'D8' '54.37'
$ws.Range('E8').Value = '  -4.82%  '
Set-TextValue 'D9' '0.372'
$ws.Range('E9').Value = '  -3.34%  '
Set-TextValue 'D10' '0.0774'
$ws.Range('E10').Value = '  -4.08%  '
$ws.Range('E11').Value = '  -5.41%  '
$ws.Range('D12').Value = '2.297.50'
$ws.Range('E12').Value = '  -2.29%  '
Set-TextValue 'D13' '13.89'
$ws.Range('E13').Value = '  -5.56%  '
Set-TextValue 'D14' '19.70'
$ws.Range('E14').Value = '  -5.35%  '
Set-TextValue 'D15' '5.21'
$ws.Range('E15').Value = '  -2.05%  '
Set-TextValue 'D16' '0.731'
$ws.Range('E16').Value = '  -3.65%  '
$ws.Range('D17').Value = '1.967.33'
$ws.Range('E17').Value = '  -3.94%  '
$ws.Range('D18').Value = '36.906.69'
$ws.Range('E18').Value = '  -0.91%  '
Set-TextValue 'D19' '6.22'
$ws.Range('E19').Value = '  +2.25%  '
Set-TextValue 'D20' '68.13'
$ws.Range('E20').Value = '  -2.35%  '
$ws.Range('D21').Value = '0.0₃0808'
$ws.Range('E21').Value = '  -3.15%  '
Set-TextValue 'D22' '221.35'
$ws.Range('E22').Value = '  -2.27%  '
Set-TextValue 'D23' '1.00'
$ws.Range('E23').Value = '  +0.04%  '
Set-TextValue 'D24' '2.42'
$ws.Range('E24').Value = '  +1.42%  '
$ws.Range('E25').Value = '  -6.39%  '
Set-TextValue 'D26' '163.85'
$ws.Range('E26').Value = '  -2.91%  '
Set-TextValue 'D27' '8.96'
$ws.Range('E27').Value = '  -6.56%  '
$ws.Range('E28').Value = '  -4.24%  '
Set-TextValue 'D29' '18.44'
$ws.Range('E29').Value = '  -3.05%  '
$ws.Range('E30').Value = '  -6.43%  '
$ws.Range('E31').Value = '  -2.28%  '
$ws.Range('E32').Value = '  -2.67%  '
$ws.Range('E33').Value = '  -3.25%  '
Set-TextValue 'D34' '4.42'
$ws.Range('E34').Value = '  -4.08%  '
$ws.Range('E35').Value = '  +2.28%  '
$ws.Range('E36').Value = '  -5.38%  '
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('E38').Value = '  -4.83%  '
Set-TextValue 'D39' '5.31'
$ws.Range('E39').Value = '  -1.48%  '
$ws.Range('D40').Value = '1.451.04'
$ws.Range('E40').Value = '  -3.21%  '
$ws.Range('E41').Value = '  -5.34%  '
Set-TextValue 'D42' '94.23'
$ws.Range('E42').Value = '  -2.49%  '
$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D43' '2.76'
$ws.Range('E43').Value = '  -4.74%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D44' '0.0905'
$ws.Range('E44').Value = '  -3.84%  '
Set-TextValue 'D45' '1.11'
$ws.Range('E45').Value = '  -4.80%  '
Set-TextValue 'D46' '15.78'
$ws.Range('E46').Value = '  -8.29%  '
$ws.Range('E47').Value = '  -0.89%  '
Set-TextValue 'D48' '0.991'
$ws.Range('E48').Value = '  -3.14%  '
$ws.Range('E49').Value = '  -0.82%  '
$ws.Range('D50').Value = '2.188.04'
$ws.Range('E50').Value = '  -2.26%  '
$ws.Range('E51').Value = '  -10.34%  '
